$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Row 1 Col 1: "66÷9=7, 3" -> "58÷2=29, 0"
$cell = $tbl.Cell(1, 1)
$cell.Range.Text = "58÷2=29, 0"

# Row 1 Col 2: "28÷4=7, 0" -> "95÷4=23, 3"
$cell = $tbl.Cell(1, 2)
$cell.Range.Text = "95÷4=23, 3"

# Row 1 Col 3: "53÷6=8, 5" -> "72÷7=10, 2"
$cell = $tbl.Cell(1, 3)
$cell.Range.Text = "72÷7=10, 2"

# Row 1 Col 4: "80÷3=26, 2" -> "62÷8=7, 6"
$cell = $tbl.Cell(1, 4)
$cell.Range.Text = "62÷8=7, 6"

# Row 1 Col 5: "56÷2=28, 0" -> "36÷4=9, 0"
$cell = $tbl.Cell(1, 5)
$cell.Range.Text = "36÷4=9, 0"

# Row 5 Col 1: "19÷5=3, 4" -> "36÷3=12, 0"
$cell = $tbl.Cell(5, 1)
$cell.Range.Text = "36÷3=12, 0"

# Row 5 Col 2: "62÷9=6, 8" -> "76÷9=8, 4"
$cell = $tbl.Cell(5, 2)
$cell.Range.Text = "76÷9=8, 4"

# Row 5 Col 3: "30÷3=10, 0" -> "41÷5=8, 1"
$cell = $tbl.Cell(5, 3)
$cell.Range.Text = "41÷5=8, 1"

# Row 5 Col 4: "82÷8=10, 2" -> "17÷9=1, 8"
$cell = $tbl.Cell(5, 4)
$cell.Range.Text = "17÷9=1, 8"

# Row 5 Col 5: "78÷4=19, 2" -> "81÷6=13, 3"
$cell = $tbl.Cell(5, 5)
$cell.Range.Text = "81÷6=13, 3"

# Row 9 Col 1: "97÷5=19, 2" -> "34÷4=8, 2"
$cell = $tbl.Cell(9, 1)
$cell.Range.Text = "34÷4=8, 2"

# Row 9 Col 2: "35÷4=8, 3" -> "82÷8=10, 2"
$cell = $tbl.Cell(9, 2)
$cell.Range.Text = "82÷8=10, 2"

# Row 9 Col 3: "80÷7=11, 3" -> "13÷4=3, 1"
$cell = $tbl.Cell(9, 3)
$cell.Range.Text = "13÷4=3, 1"

# Row 9 Col 4: "81÷8=10, 1" -> "89÷4=22, 1"
$cell = $tbl.Cell(9, 4)
$cell.Range.Text = "89÷4=22, 1"

# Row 9 Col 5: "77÷5=15, 2" -> "49÷2=24, 1"
$cell = $tbl.Cell(9, 5)
$cell.Range.Text = "49÷2=24, 1"

# Row 13 Col 1: "56÷7=8, 0" -> "93÷2=46, 1"
$cell = $tbl.Cell(13, 1)
$cell.Range.Text = "93÷2=46, 1"

# Row 13 Col 2: "73÷5=14, 3" -> "35÷7=5, 0"
$cell = $tbl.Cell(13, 2)
$cell.Range.Text = "35÷7=5, 0"

# Row 13 Col 3: "28÷2=14, 0" -> "25÷6=4, 1"
$cell = $tbl.Cell(13, 3)
$cell.Range.Text = "25÷6=4, 1"

# Row 13 Col 4: "32÷3=10, 2" -> "17÷4=4, 1"
$cell = $tbl.Cell(13, 4)
$cell.Range.Text = "17÷4=4, 1"

# Row 13 Col 5: "13÷4=3, 1" -> "93÷9=10, 3"
$cell = $tbl.Cell(13, 5)
$cell.Range.Text = "93÷9=10, 3"

# Row 17 Col 1: "70÷6=11, 4" -> "96÷6=16, 0"
$cell = $tbl.Cell(17, 1)
$cell.Range.Text = "96÷6=16, 0"

# Row 17 Col 2: "16÷7=2, 2" -> "31÷4=7, 3"
$cell = $tbl.Cell(17, 2)
$cell.Range.Text = "31÷4=7, 3"

# Row 17 Col 3: "34÷4=8, 2" -> "47÷5=9, 2"
$cell = $tbl.Cell(17, 3)
$cell.Range.Text = "47÷5=9, 2"

# Row 17 Col 4: "48÷2=24, 0" -> "42÷6=7, 0"
$cell = $tbl.Cell(17, 4)
$cell.Range.Text = "42÷6=7, 0"

# Row 17 Col 5: "43÷9=4, 7" -> "19÷3=6, 1"
$cell = $tbl.Cell(17, 5)
$cell.Range.Text = "19÷3=6, 1"
